$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7613898750257053
$ws.Range("C2").Value = 0.7173457508731083
$ws.Range("D2").Value = 0.7375309023440424
$ws.Range("E2").Value = 0.5034286756980586
$ws.Range("F2").Value = 0.5079934600998902
